# Add team record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style used by the other headers
# (e.g. AC1) by copying its format onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2:55 -> every row gets the same team record: 101 wins, 61 losses, 0 ties
$ws.Range("AD2:AD55").Value = 101
$ws.Range("AE2:AE55").Value = 61
$ws.Range("AF2:AF55").Value = 0
